# The deck currently has its live/applied theme ("Integral" colours) stored
# in ppt/theme/theme2.xml (the file referenced by the slide master and by
# the package's default theme relationship), while ppt/theme/theme1.xml
# holds the stock "Office Theme" colours and is only referenced by the
# (visually inert) notes master.
#
# The target edit swaps the two theme bodies: the colours actually applied
# to the presentation become the default "Office Theme" palette again,
# while the old "Integral" palette is pushed into the unused theme slot.
#
# PowerPoint's object model does not expose raw OOXML part plumbing, so we
# reproduce the visible/semantic effect of that swap with the supported
# Theme-colour API: we recolour the live theme colour scheme (reachable
# through any slide's ThemeColorScheme, which maps 1:1 onto the master's
# <a:clrScheme> slots dk1,lt1,dk2,lt2,accent1-6,hlink,folHlink) so that it
# matches the standard Office Theme palette again.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Index -> (name, Office Theme RGB as 0xBBGGRR for the COM RGB() encoding)
$officeColors = @{
    1  = 0            # dk1      000000
    2  = 16777215      # lt1      FFFFFF
    3  = 6968388        # dk2      44546A
    4  = 15132391       # lt2      E7E6E6
    5  = 13998939       # accent1  5B9BD5
    6  = 3243501         # accent2  ED7D31
    7  = 10855845        # accent3  A5A5A5
    8  = 49407           # accent4  FFC000
    9  = 12874308        # accent5  4472C4
    10 = 4697456         # accent6  70AD47
    11 = 12673797        # hlink    0563C1
    12 = 7491477          # folHlink 954F72
}

for ($i = 1; $i -le 12; $i++) {
    $color = $tcs.Colors($i)
    $color.RGB = $officeColors[$i]
}
